# Append a new statistics row (row 26) to the "統計" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("統計")

$newRow = 26

$ws.Cells.Item($newRow, 1).Value = "2025-08-30T06:29:34.603425"
$ws.Cells.Item($newRow, 2).Value = 6
$ws.Cells.Item($newRow, 3).Value = "全案件リスト"
$ws.Cells.Item($newRow, 4).Value = 66.7
$ws.Cells.Item($newRow, 5).Value = 2
$ws.Cells.Item($newRow, 6).Value = 3
$ws.Cells.Item($newRow, 7).Value = 6
